# Weekly update: add two new "Coco" price records for Terminal La Palmera de La Serena.
# This mirrors the real-world process of prepending newly scraped weekly rows to the
# top of the data block (after the header) and pushing the existing history down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Common (unchanging) field values for every row in this data set ---
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100108
$producto   = "Tropicales y subtropicales"
$categoriaId = 100108007
$categoria   = "Coco"
$variedad    = "Sin especificar"
$calidad     = "Primera"
$unidad      = "$/malla 20 unidades"
$origen      = "Perú"
$kgUnidad    = 20

function Set-CocoRow($rowNum, $fecha, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($rowNum, 1).Value  = $mercadoId
    $ws.Cells.Item($rowNum, 2).Value  = $mercado
    $ws.Cells.Item($rowNum, 3).Value  = $region
    $ws.Cells.Item($rowNum, 4).Value  = $fecha
    $ws.Cells.Item($rowNum, 5).Value  = $codreg
    $ws.Cells.Item($rowNum, 6).Value  = $tipo
    $ws.Cells.Item($rowNum, 7).Value  = $productoId
    $ws.Cells.Item($rowNum, 8).Value  = $producto
    $ws.Cells.Item($rowNum, 9).Value  = $categoriaId
    $ws.Cells.Item($rowNum, 10).Value = $categoria
    $ws.Cells.Item($rowNum, 11).Value = $variedad
    $ws.Cells.Item($rowNum, 12).Value = $calidad
    $ws.Cells.Item($rowNum, 13).Value = $volumen
    $ws.Cells.Item($rowNum, 14).Value = $precioMin
    $ws.Cells.Item($rowNum, 15).Value = $precioMax
    $ws.Cells.Item($rowNum, 16).Value = $precioProm
    $ws.Cells.Item($rowNum, 17).Value = $unidad
    $ws.Cells.Item($rowNum, 18).Value = $origen
    $ws.Cells.Item($rowNum, 19).Value = $precioKg
    $ws.Cells.Item($rowNum, 20).Value = $kgUnidad
}

# Insert first new row at row 26, pushing the old row 26 (and everything below) down by one.
$ws.Rows.Item(26).Insert()
Set-CocoRow 26 44874 240 29000 30000 29500 1475

# Insert second new row at row 32 (which, after the first insert, sits right after the old
# row 30 / before the old row 31), pushing everything from there down by one more.
$ws.Rows.Item(32).Insert()
Set-CocoRow 32 44879 100 28000 30000 29000 1450
